# Apply the "0 Information Meeting.pptx" edits:
#  1. Slide 2 ("Number of students and workload"): update the student
#     headcount from 126 to 130 and merge the sentence into a single run.
#  2. Slide 6 ("bachelor project proposals"): tweak the "contact person"
#     bullet wording.
#  3. Slide 6: trim "Groups with 4 or more persons are not allowed (by the
#     formal rules)" down to "...are not allowed", splitting the sentence
#     into two runs ("...are " / "not allowed").

$p = $ppt.ActivePresentation

# --- Change 1: Slide 2, "Content Placeholder 2" shape, first paragraph ---
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$para1 = $tr2.Paragraphs(1, 1)
$fullRange1 = $para1.Characters(1, $para1.Length)
$fullRange1.Text = "There will be a total of 130 students doing their cs / it bachelor project in the Spring of 2026"

# --- Change 2 & 3: Slide 6, "Rectangle 3" shape ---
$slide6 = $p.Slides.Item(6)
$shape6 = $slide6.Shapes.Item(2)
$tr6 = $shape6.TextFrame.TextRange

# Change 2: paragraph 2 - "contact person" wording tweak
$para2 = $tr6.Paragraphs(2, 1)
$fullRange2 = $para2.Characters(1, $para2.Length)
$fullRange2.Text = "You are encouraged to speak with the contact person(s) for the corresponding section to obtain additional information"

# Change 3: paragraph 11 - drop "(by the formal rules)"
$para11 = $tr6.Paragraphs(11, 1)
# Re-assigning the leading chunk's text forces a run split right after it,
# leaving "Groups with 4 or more persons are " / "not allowed (by the formal rules)"
$lead = $para11.Characters(1, 34)
$lead.Text = "Groups with 4 or more persons are "
# Now drop the trailing " (by the formal rules)" text.
$full11 = $para11.Text
$tailIdx = $full11.IndexOf(" (by the formal rules)")
$tail = $para11.Characters($tailIdx + 1, $full11.Length - $tailIdx)
$tail.Text = ""
